# Regenerate the "K" column (column G) values for rows 2-45 on the
# active worksheet. These values were recomputed upstream (switching
# the stat source from "Strike#" to "K"), so we just write the new
# literal values into column G for each game row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New K values for rows 2..45 (game index 0..43), in order.
$kValues = @(0,0,0,1,1,2,1,1,0,2,2,0,1,1,0,1,1,0,0,0,1,0,0,0,1,0,0,0,1,1,1,0,2,1,0,0,0,1,1,1,0,0,1,1)

$row = 2
foreach ($val in $kValues) {
    $ws.Cells.Item($row, 7).Value = $val   # Column G is the 7th column
    $row++
}
